# Applies the "DescriptorSet Issue Finally Solved" edit:
# Adds two paragraphs after the last paragraph of the document:
#   1. A bold line reporting the vkAllocateDescriptorSet error code.
#   2. A bulleted list item explaining the root cause of the error.

$d = $word.ActiveDocument
$wdNamespace = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# -- Paragraph 1: bold status line ------------------------------------------------
$p1Xml = '<w:p ' + $wdNamespace + '>' + `
    '<w:pPr><w:spacing w:after="0"/><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>vkAllocateDescriptorSet</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">error code = </w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>VK_ERROR_INITIALIZATION_FAILED</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '</w:p>'

# -- Paragraph 2: bulleted explanation --------------------------------------------
$rsquo = [char]0x2019
$p2Xml = '<w:p ' + $wdNamespace + '>' + `
    '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>The error come from the size of the descriptorPool</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> , there is 2 parameter </w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>the size of the pool and the max size</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, there is no error if you have the size&gt;maxsize but when you allocateDescriptor an error occur, and you don' + $rsquo + 't get easily why</w:t></w:r>' + `
    '</w:p>'

# -- Insert both paragraphs after the current last paragraph ---------------------
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertXML($p1Xml) | Out-Null

$tail2 = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$tail2.Collapse(0)
$tail2.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertXML($p2Xml) | Out-Null
